$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A164").Value = "baguss"
$ws.Range("B164").Value = "bagus"
$ws.Range("A165").Value = "buriq"
$ws.Range("B165").Value = "burik"

$ws.Range("G172").Select()
